$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("September")

# New row for the next movable date (21.09.2018, Fri), mirroring the
# pattern used by the existing "19.09.2018, Wed" row (row 20).
$ws.Range("A22").Value = "21.09.2018, Fri"
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = "`n"

# Setting a newline-only value auto-expands the row height; restore the
# natural/default height so row 22 matches the other data rows.
$ws.Rows.Item(22).AutoFit()
